$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row for new columns G:J ---
$ws.Cells.Item(1, 7).Value = "Age"
$ws.Cells.Item(1, 8).Value = "Sex"
$ws.Cells.Item(1, 9).Value = "Batch"
$ws.Cells.Item(1, 10).Value = "Weight"
$ws.Range("G1:J1").Font.Bold = $true

# --- Per-row data: row, Age, Sex, Batch, Weight ---
$rows = @(
  @(2, 8, "Gelding", 1, 500),
  @(3, 8, "Gelding", 1, 500),
  @(4, 8, "Gelding", 1, 500),
  @(5, 12, "Mare", 1, 550),
  @(6, 12, "Mare", 1, 550),
  @(7, 12, "Mare", 1, 550),
  @(8, 5, "Mare", 1, 410),
  @(9, 5, "Mare", 1, 410),
  @(10, 5, "Mare", 1, 410),
  @(11, 7, "Mare", 1, 505),
  @(12, 7, "Mare", 1, 505),
  @(13, 7, "Mare", 1, 505),
  @(14, 6, "Gelding", 2, 550),
  @(15, 6, "Gelding", 2, 550),
  @(16, 6, "Gelding", 2, 550),
  @(17, 6, "Mare", 2, 545),
  @(18, 6, "Mare", 2, 545),
  @(19, 6, "Mare", 2, 545),
  @(20, 5, "Mare", 2, 450),
  @(21, 5, "Mare", 2, 450),
  @(22, 5, "Mare", 2, 450),
  @(23, 5, "Gelding", 2, 460),
  @(24, 5, "Gelding", 2, 460),
  @(25, 5, "Gelding", 2, 460),
  @(26, 6, "Gelding", 2, 450),
  @(27, 6, "Gelding", 2, 450),
  @(28, 6, "Gelding", 2, 450),
  @(29, 10, "Gelding", 2, 480),
  @(30, 10, "Gelding", 2, 480),
  @(31, 10, "Gelding", 2, 480),
  @(32, 9, "Mare", 3, 490),
  @(33, 9, "Mare", 3, 490),
  @(34, 9, "Mare", 3, 490),
  @(35, 8, "Mare", 3, 460),
  @(36, 8, "Mare", 3, 460),
  @(37, 8, "Mare", 3, 460),
  @(38, 7, "Mare", 3, 450),
  @(39, 7, "Mare", 3, 450),
  @(40, 7, "Mare", 3, 450),
  @(41, 6, "Gelding", 3, 500),
  @(42, 6, "Gelding", 3, 500),
  @(43, 6, "Gelding", 3, 500),
  @(44, 7, "Mare", 3, 575),
  @(45, 7, "Mare", 3, 575),
  @(46, 7, "Mare", 3, 575),
  @(47, 8, "Gelding", 3, 500),
  @(48, 8, "Gelding", 3, 500),
  @(49, 8, "Gelding", 3, 500),
  @(50, 6, "Mare", 3, 525),
  @(51, 6, "Mare", 3, 525),
  @(52, 6, "Mare", 3, 525),
  @(53, 6, "Mare", 4, 500),
  @(54, 6, "Mare", 4, 500),
  @(55, 6, "Mare", 4, 500),
  @(56, 5, "Mare", 4, 520),
  @(57, 5, "Mare", 4, 520),
  @(58, 5, "Mare", 4, 520),
  @(59, 6, "Mare", 4, 420),
  @(60, 6, "Mare", 4, 420),
  @(61, 6, "Mare", 4, 420),
  @(62, 5, "Mare", 4, 435),
  @(63, 5, "Mare", 4, 435),
  @(64, 5, "Mare", 4, 435),
  @(65, 8, "Mare", 4, 460),
  @(66, 8, "Mare", 4, 460),
  @(67, 8, "Mare", 4, 460),
  @(68, 8, "Mare", 4, 530),
  @(69, 8, "Mare", 4, 530),
  @(70, 8, "Mare", 4, 530),
  @(71, 4, "Mare", 4, 445),
  @(72, 4, "Mare", 4, 445),
  @(73, 4, "Mare", 4, 445),

)

foreach ($item in $rows) {
    $r = $item[0]
    $ws.Cells.Item($r, 7).Value = $item[1]
    $ws.Cells.Item($r, 8).Value = $item[2]
    $ws.Cells.Item($r, 9).Value = $item[3]
    $ws.Cells.Item($r, 10).Value = $item[4]
}

# --- Column widths for new columns ---
$ws.Columns.Item(7).ColumnWidth = 13.833333333333332
$ws.Columns.Item(8).ColumnWidth = 14.333333333333332
$ws.Columns.Item(9).ColumnWidth = 17.833333333333336

# --- AutoFilter over the full table range ---
$ws.Range("A1:J1").AutoFilter() | Out-Null
$ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$J`$1")
$n = $ws.Names.Item(1)
$n.Visible = $false

# --- Selection matches the saved workbook state ---
$ws.Range("O9").Select()
